$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly date column (Thursday), one week after column E (45701 -> 45708)
$ws.Range("F1").Value = 45708
$ws.Range("F1").NumberFormat = $ws.Range("E1").NumberFormat

# Header "Qui" (same text/style as C2/D2/E2, unformatted relative to B2)
$ws.Range("F2").Value = "Qui"

# Attendance marks "P" for each student row
for ($r = 3; $r -le 20; $r++) {
    $ws.Cells.Item($r, 6).Value = "P"
}

# Update selection to reflect the newly active cell
$ws.Range("F3").Select()
